# Insert two new weekly price rows for Femacal de La Calera - Mango,
# right before the existing row 408, shifting rows 408:439 down to 410:441.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new blank rows at position 408 (pushes old 408:439 -> 410:441)
$ws.Rows.Item(408).Resize(2).Insert()

# Row 408: new Primera entry, Brasil origin
$ws.Range("A408").Value = 3
$ws.Range("B408").Value = "Femacal de La Calera"
$ws.Range("C408").Value = "Coquimbo"
$ws.Range("D408").Value = 44769
$ws.Range("E408").Value = 5
$ws.Range("F408").Value = "Fruta"
$ws.Range("G408").Value = 100108
$ws.Range("H408").Value = "Tropicales y subtropicales"
$ws.Range("I408").Value = 100108002
$ws.Range("J408").Value = "Mango"
$ws.Range("K408").Value = "Sin especificar"
$ws.Range("L408").Value = "Primera"
$ws.Range("M408").Value = 228
$ws.Range("N408").Value = 9000
$ws.Range("O408").Value = 9000
$ws.Range("P408").Value = 9000
$ws.Range("Q408").Value = "`$/bandeja 4 kilos"
$ws.Range("R408").Value = "Brasil"
$ws.Range("S408").Value = 2250
$ws.Range("T408").Value = 4

# Row 409: new Segunda entry, Brasil origin
$ws.Range("A409").Value = 3
$ws.Range("B409").Value = "Femacal de La Calera"
$ws.Range("C409").Value = "Coquimbo"
$ws.Range("D409").Value = 44769
$ws.Range("E409").Value = 5
$ws.Range("F409").Value = "Fruta"
$ws.Range("G409").Value = 100108
$ws.Range("H409").Value = "Tropicales y subtropicales"
$ws.Range("I409").Value = 100108002
$ws.Range("J409").Value = "Mango"
$ws.Range("K409").Value = "Sin especificar"
$ws.Range("L409").Value = "Segunda"
$ws.Range("M409").Value = 228
$ws.Range("N409").Value = 9000
$ws.Range("O409").Value = 9000
$ws.Range("P409").Value = 9000
$ws.Range("Q409").Value = "`$/bandeja 4 kilos"
$ws.Range("R409").Value = "Brasil"
$ws.Range("S409").Value = 2250
$ws.Range("T409").Value = 4

"Inserted rows 408-409; dimension should now be A1:T441"
